{"js": "// New translations for \"04_Playful mathematicians - subtitles (format and\n// timing corrected).docx\" (Swahili, Tanzania): replace the Swahili caption\n// text of each subtitle cue with its English translation, leaving the\n// surrounding timing/formatting untouched.\nconst pairs = [\n  [\"Wanahisabati wanaocheza:\", \"The playful mathematicians:\"],\n  [\"** mazungumzo huanza saa 47 ya pili, kwa hivyo niliongeza sekunde 28 kwa nyakati zote kama zilivyokuwa. -John Argentino\", \"** the dialogue starts at second 47, so I added 28 seconds to all the times as they were. -John Argentino\"],\n  [\"[Muziki]\", \"[Music]\"],\n  [\"kuna wanahisabati wawili, tupige simu\", \"there are two mathematicians, let's call\"],\n  [\"Fil na Mike wanaokutana\", \"them Fil and Mike who meet each other\"],\n  [\"tena baada ya muda mrefu. Baada ya baadhi\", \"again after a long time. After some\"],\n  [\"kuzungumza, Phil anasema ana watoto watatu, basi\", \"chatting, Phil says he has three children, then\"],\n  [\"Kwa mshangao, Mike anauliza: 'Wana umri gani?' Fil,\", \"Mike, astonished, asks: 'How old are they?' Fil,\"],\n  [\"kuwa mwanahisabati mchezaji, anajibu\", \"being a playful mathematician, answers\"],\n  [\"'Wewe niambie! Nitakupa kidokezo: ikiwa wewe\", \"'You tell me! I'll give you a hint: if you\"],\n  [\"zidisheni enzi tatu pamoja ninyi\", \"multiply the three ages together you\"],\n  [\"pata 36.' Mike huchukua wakati mwingine kufikiria\", \"get 36.' Mike takes sometimes to think\"],\n  [\"na kusema: 'Samahani Fil, lakini nahitaji\", \"and says: 'I'm sorry Fil, but I do need\"],\n  [\"kidokezo kingine. Kwa hivyo Fil anamwambia Mike:\", \"another hint. So Fil tells Mike:\"],\n  [\"'Ndiyo, hakika, hapa ni: kama alikuwa na hadi\", \"'Yes, sure, here it is: if you had up to\"],\n  [\"miaka mitatu unapata idadi ya hesabu\", \"three ages you get the number of math\"],\n  [\"karatasi tunachapisha pamoja. Je, unaikumbuka?'\", \"papers we publish together. Do you remember it?'\"],\n  [\"'Ndio nakumbuka wangapi, lakini bado\", \"'Yes I do remember How many, but still\"],\n  [\"Sina taarifa za kutosha! nahitaji\", \"I do not have enough information! I need\"],\n  [\"angalau moja zaidi.' Fil anasema: 'Ndiyo usifanye hivyo\", \"at least one more.' Fil says: 'Yes don't\"],\n  [\"wasiwasi lakini hii ni ya mwisho:\", \"worry but this is the last one:\"],\n  [\"Mdogo ana macho ya blues.' Na\", \"The youngest one has blues eyes.' And\"],\n  [\"ghafla Mike anapata jibu. Wewe\", \"suddenly Mike gets the answer. You\"],\n  [\"sikia mazungumzo lakini hujui\", \"hear the conversation but you don't know\"],\n  [\"ni karatasi ngapi walichapisha pamoja.\", \"how many papers they published together.\"],\n  [\"Hata hivyo, unataka kujua umri wa\", \"However, you do want to know the ages of\"],\n  [\"watoto watatu. Je, unaweza kuwahesabu\", \"the three children. Can you figure them\"],\n  [\"nje?\", \"out?\"],\n];\n\nconst body = context.document.body;\nfor (const [from, to] of pairs) {\n  const results = body.search(from, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (const item of results.items) {\n    item.insertText(to, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# New translations for \"04_Playful mathematicians - subtitles (format and\n# timing corrected).docx\" (Swahili, Tanzania): replace the Swahili caption\n# text of each subtitle cue with its English translation, leaving the\n# surrounding timing/formatting untouched.\n\n# Turn off \"smart quotes\" autocorrect so the straight apostrophes in the\n# translated text are kept as authored instead of being curled.\n$word.Options.AutoFormatAsYouTypeReplaceQuotes = $false\n$word.Options.AutoFormatReplaceQuotes = $false\n\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.Text = \"Wanahisabati wanaocheza:\"\n$find.Forward = $true\n$find.Wrap = 0\nwhile ($find.Execute()) {\n    $find.Parent.Text = \"The playful mathematicians:\"\n}\n\n$find = $d.Content.Find\n$find.Text = \"** mazungumzo huanza saa 47 ya pili, kwa hivyo niliongeza sekunde 28 kwa nyakati zote kama zilivyokuwa. -John Argentino\"\n$find.Forward = $true\n$find.Wrap = 0\nwhile ($find.Execute()) {\n    $find.Parent.Text = \"** the dialogue starts at second 47, so I added 28 seconds to all the times as they were. -John Argentino\"\n}\n\n$find = $d.Content.Find\n$find.Text = \"[Muziki]\"\n$find.Forward = $true\n$find.Wrap = 0\nwhile ($find.Execute()) {\n    $find.Parent.Text = \"[Music]\"\n}\n\n$find = $d.Content.Find\n$find.Text = \"kuna wanahisabati wawili, tupige simu\"\n$find.Forward = $true\n$find.Wrap = 0\nwhile ($find.Execute()) {\n    $find.Parent.Text = \"there are two mathematicians, let's call\"\n}\n\n$find = $d.Content.Find\n$find.Text = \"Fil na Mike wanaokutana\"\n$find.Forward = $true\n$find.Wrap = 0\nwhile ($find.Execute()) {\n    $find.Parent.Text = \"them Fil and Mike who meet each other\"\n}\n\n$find = $d.Content.Find\n$find.Text = \"tena baada ya muda mrefu. Baada ya baadhi\"\n$find.Forward = $true\n$find.Wrap = 0\nwhile ($find.Execute()) {\n    $find.Parent.Text = \"again after a long time. After some\"\n}\n\n$find = $d.Content.Find\n$find.Text = \"kuzungumza, Phil anasema ana watoto watatu, basi\"\n$find.Forward = $true\n$find.Wrap = 0\nwhile ($find.Execute()) {\n    $find.Parent.Text = \"chatting, Phil says he has three children, then\"\n}\n\n$find = $d.Content.Find\n$find.Text = \"Kwa mshangao, Mike anauliza: 'Wana umri gani?' Fil,\"\n$find.Forward = $true\n$find.Wrap = 0\nwhile ($find.Execute()) {\n    $find.Parent.Text = \"Mike, astonished, asks: 'How old are they?' Fil,\"\n}\n\n$find = $d.Content.Find\n$find.Text = \"kuwa mwanahisabati mchezaji, anajibu\"\n$find.Forward = $true\n$find.Wrap = 0\nwhile ($find.Execute()) {\n    $find.Parent.Text = \"being a playful mathematician, answers\"\n}\n\n$find = $d.Content.Find\n$find.Text = \"'Wewe niambie! Nitakupa kidokezo: ikiwa wewe\"\n$find.Forward = $true\n$find.Wrap = 0\nwhile ($find.Execute()) {\n    $find.Parent.Text = \"'You tell me! I'll give you a hint: if you\"\n}\n\n$find = $d.Content.Find\n$find.Text = \"zidisheni enzi tatu pamoja ninyi\"\n$find.Forward = $true\n$find.Wrap = 0\nwhile ($find.Execute()) {\n    $find.Parent.Text = \"multiply the three ages together you\"\n}\n\n$find = $d.Content.Find\n$find.Text = \"pata 36.' Mike huchukua wakati mwingine kufikiria\"\n$find.Forward = $true\n$find.Wrap = 0\nwhile ($find.Execute()) {\n    $find.Parent.Text = \"get 36.' Mike takes sometimes to think\"\n}\n\n$find = $d.Content.Find\n$find.Text = \"na kusema: 'Samahani Fil, lakini nahitaji\"\n$find.Forward = $true\n$find.Wrap = 0\nwhile ($find.Execute()) {\n    $find.Parent.Text = \"and says: 'I'm sorry Fil, but I do need\"\n}\n\n$find = $d.Content.Find\n$find.Text = \"kidokezo kingine. Kwa hivyo Fil anamwambia Mike:\"\n$find.Forward = $true\n$find.Wrap = 0\nwhile ($find.Execute()) {\n    $find.Parent.Text = \"another hint. So Fil tells Mike:\"\n}\n\n$find = $d.Content.Find\n$find.Text = \"'Ndiyo, hakika, hapa ni: kama alikuwa na hadi\"\n$find.Forward = $true\n$find.Wrap = 0\nwhile ($find.Execute()) {\n    $find.Parent.Text = \"'Yes, sure, here it is: if you had up to\"\n}\n\n$find = $d.Content.Find\n$find.Text = \"miaka mitatu unapata idadi ya hesabu\"\n$find.Forward = $true\n$find.Wrap = 0\nwhile ($find.Execute()) {\n    $find.Parent.Text = \"three ages you get the number of math\"\n}\n\n$find = $d.Content.Find\n$find.Text = \"karatasi tunachapisha pamoja. Je, unaikumbuka?'\"\n$find.Forward = $true\n$find.Wrap = 0\nwhile ($find.Execute()) {\n    $find.Parent.Text = \"papers we publish together. Do you remember it?'\"\n}\n\n$find = $d.Content.Find\n$find.Text = \"'Ndio nakumbuka wangapi, lakini bado\"\n$find.Forward = $true\n$find.Wrap = 0\nwhile ($find.Execute()) {\n    $find.Parent.Text = \"'Yes I do remember How many, but still\"\n}\n\n$find = $d.Content.Find\n$find.Text = \"Sina taarifa za kutosha! nahitaji\"\n$find.Forward = $true\n$find.Wrap = 0\nwhile ($find.Execute()) {\n    $find.Parent.Text = \"I do not have enough information! I need\"\n}\n\n$find = $d.Content.Find\n$find.Text = \"angalau moja zaidi.' Fil anasema: 'Ndiyo usifanye hivyo\"\n$find.Forward = $true\n$find.Wrap = 0\nwhile ($find.Execute()) {\n    $find.Parent.Text = \"at least one more.' Fil says: 'Yes don't\"\n}\n\n$find = $d.Content.Find\n$find.Text = \"wasiwasi lakini hii ni ya mwisho:\"\n$find.Forward = $true\n$find.Wrap = 0\nwhile ($find.Execute()) {\n    $find.Parent.Text = \"worry but this is the last one:\"\n}\n\n$find = $d.Content.Find\n$find.Text = \"Mdogo ana macho ya blues.' Na\"\n$find.Forward = $true\n$find.Wrap = 0\nwhile ($find.Execute()) {\n    $find.Parent.Text = \"The youngest one has blues eyes.' And\"\n}\n\n$find = $d.Content.Find\n$find.Text = \"ghafla Mike anapata jibu. Wewe\"\n$find.Forward = $true\n$find.Wrap = 0\nwhile ($find.Execute()) {\n    $find.Parent.Text = \"suddenly Mike gets the answer. You\"\n}\n\n$find = $d.Content.Find\n$find.Text = \"sikia mazungumzo lakini hujui\"\n$find.Forward = $true\n$find.Wrap = 0\nwhile ($find.Execute()) {\n    $find.Parent.Text = \"hear the conversation but you don't know\"\n}\n\n$find = $d.Content.Find\n$find.Text = \"ni karatasi ngapi walichapisha pamoja.\"\n$find.Forward = $true\n$find.Wrap = 0\nwhile ($find.Execute()) {\n    $find.Parent.Text = \"how many papers they published together.\"\n}\n\n$find = $d.Content.Find\n$find.Text = \"Hata hivyo, unataka kujua umri wa\"\n$find.Forward = $true\n$find.Wrap = 0\nwhile ($find.Execute()) {\n    $find.Parent.Text = \"However, you do want to know the ages of\"\n}\n\n$find = $d.Content.Find\n$find.Text = \"watoto watatu. Je, unaweza kuwahesabu\"\n$find.Forward = $true\n$find.Wrap = 0\nwhile ($find.Execute()) {\n    $find.Parent.Text = \"the three children. Can you figure them\"\n}\n\n$find = $d.Content.Find\n$find.Text = \"nje?\"\n$find.Forward = $true\n$find.Wrap = 0\nwhile ($find.Execute()) {\n    $find.Parent.Text = \"out?\"\n}\n\n"}
